# "from stock to counter e.g:fstc"
# Adds a new counter (khuhiLeftHani / avishekOnline) entry to both the
# 07-07-2023 and 08-07-2023 sheets, and updates the running total on
# sheet1 to include the new counter value.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 07-07-2023
$ws2 = $wb.Worksheets.Item(2)   # 08-07-2023 (active tab)

# --- Sheet1 ("07-07-2023") ---------------------------------------------
# New counter name/value pair in columns I & J
$ws1.Range("I2").Value = "khuhiLeftHani"
$ws1.Range("J2").Value = 1000

# Column I width
$ws1.Columns.Item(9).ColumnWidth = 12.166666666666666

# Running total in B2 now sums all four counters
$ws1.Range("B2").Formula = "=SUM(D2,F2,H2,J2)"

# Cursor ends up one row below the data after data entry
$null = $ws1.Range("B3").Select()

# --- Sheet2 ("08-07-2023") ----------------------------------------------
$ws2.Range("B2").Value = 21300
$ws2.Range("B2").NumberFormat = $ws2.Range("B1").NumberFormat
$ws2.Range("C2").Value = "avishekOnline"

# Column widths for the newly used columns B & C
$ws2.Columns.Item(2).ColumnWidth = 10.333333333333332
$ws2.Columns.Item(3).ColumnWidth = 11.5

# Cursor ends up one row below the data, keeping sheet2 as the active tab
$null = $ws2.Range("C3").Select()
